# Weekly update: insert a new data row for the latest week at row 117,
# pushing the existing rows (117-206) down by one (to 118-207).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 117; everything below (old 117..206) shifts to 118..207.
$ws.Rows("117:117").Insert()

# Populate the newly inserted row 117 with the new week's record (same
# market/category/quality metadata as its neighbours, new date + prices).
$ws.Range("A117").Value = 8
$ws.Range("B117").Value = "Terminal La Palmera de La Serena"
$ws.Range("C117").Value = "Coquimbo"
$ws.Range("D117").Value = 44978
$ws.Range("D117").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E117").Value = 4
$ws.Range("F117").Value = 100112040
$ws.Range("G117").Value = "Cilantro"
$ws.Range("H117").Value = "Sin especificar"
$ws.Range("I117").Value = "Primera"
$ws.Range("J117").Value = 2200
$ws.Range("K117").Value = 2300
$ws.Range("L117").Value = 2500
$ws.Range("M117").Value = 2400
$ws.Range("N117").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O117").Value = "Provincia del Elquí"
$ws.Range("P117").Value = 1600
$ws.Range("Q117").Value = 1.5
$ws.Range("R117").Value = "Hortaliza"
